# Apply the diff: split the sentence ending "...SL Domain-Specific
# Language. We" so that "Language" is followed by " for solving stencils"
# before the closing period, i.e.:
#   "This feature is implemented as a template for use in the SL
#    Domain-Specific Language. We"
# becomes:
#   "This feature is implemented as a template for use in the SL
#    Domain-Specific Language for solving stencils. We"

$d = $word.ActiveDocument

$old = "This feature is implemented as a template for use in the SL Domain-Specific Language. We"
$new = "This feature is implemented as a template for use in the SL Domain-Specific Language for solving stencils. We"

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the target sentence to update."
}
